$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the match-record contents (columns B, C, E..AD) between paired rows.
# Column A (running id) and column D (match date) stay fixed in place.
$rowPairs = @(
    @(17, 18),
    @(81, 82),
    @(105, 106),
    @(107, 108),
    @(121, 122),
    @(135, 136),
    @(161, 162)
)

# Columns to swap: B..C (2..3) and E..AD (5..30); D (4) = date is left untouched.
$swapCols = @(2, 3) + @(5..30)

foreach ($pair in $rowPairs) {
    $rowA = $pair[0]
    $rowB = $pair[1]
    foreach ($col in $swapCols) {
        $cellA = $ws.Cells.Item($rowA, $col)
        $cellB = $ws.Cells.Item($rowB, $col)
        $valA = $cellA.Value2
        $valB = $cellB.Value2
        $cellA.Value = $valB
        $cellB.Value = $valA
    }
}
